# Added common test data for Belgium and Germany market
# Insert a new "FAT-S" row into the Accessories list on both the
# "Accessories" and "Accessories_215 Panel" sheets, just above the
# trailing "Wg" / "Accessories" sentinel rows, then update the
# active-sheet/selection state to match.

$wb = $excel.ActiveWorkbook

# --- "Accessories" sheet (row 14 -> new "FAT-S" row; old 14/15 shift to 15/16) ---
$ws8 = $wb.Worksheets.Item("Accessories")
$ws8.Rows.Item(14).Insert()
$ws8.Cells.Item(15, 1).Copy()
$ws8.Cells.Item(14, 1).PasteSpecial(-4122)
$ws8.Cells.Item(14, 1).Value = "FAT-S"

# --- "Accessories_215 Panel" sheet (row 11 -> new "FAT-S" row; old 11/12 shift to 12/13) ---
$ws9 = $wb.Worksheets.Item("Accessories_215 Panel")
$ws9.Rows.Item(11).Insert()
$ws9.Cells.Item(12, 1).Copy()
$ws9.Cells.Item(11, 1).PasteSpecial(-4122)
$ws9.Cells.Item(11, 1).Value = "FAT-S"

# --- Update selection/active-sheet state ---
# "Accessories_215 Panel" selection moves to A9 (no longer the active tab).
$ws9.Activate()
$ws9.Range("A9").Select()

# "Accessories" becomes the active tab, selection at A15.
$ws8.Activate()
$ws8.Range("A15").Select()
